$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 23.20031648293078
$ws.Range("C2").Value = 9.121142890126894
$ws.Range("D2").Value = 8.180230647678927
$ws.Range("E2").Value = 9.726312213286716
$ws.Range("F2").Value = 42.38014468077398
$ws.Range("H2").Value = 7.344005520526261
$ws.Range("L2").Value = 10.02731762449572
$ws.Range("M2").Value = 18.44414851031804
$ws.Range("N2").Value = 21.31830534023507
$ws.Range("B3").Value = 22.79640028415686
$ws.Range("C3").Value = 8.520819702750625
$ws.Range("D3").Value = 8.203359364886706
$ws.Range("E3").Value = 9.708442519933659
$ws.Range("F3").Value = 42.07396438172274
$ws.Range("H3").Value = 7.344005520526261
$ws.Range("L3").Value = 10.04054932860214
$ws.Range("M3").Value = 18.37493880585787
$ws.Range("N3").Value = 21.37021744265675
$ws.Range("B4").Value = 22.55225410524635
$ws.Range("C4").Value = 8.129543559500888
$ws.Range("D4").Value = 8.218590814447067
$ws.Range("E4").Value = 9.697215147094937
$ws.Range("F4").Value = 41.8969425217612
$ws.Range("H4").Value = 7.344005520526261
$ws.Range("L4").Value = 10.05016902832882
$ws.Range("M4").Value = 18.33678581621586
$ws.Range("N4").Value = 21.40404577617452
$ws.Range("B5").Value = 22.4538732988328
$ws.Range("C5").Value = 7.964377621417094
$ws.Range("D5").Value = 8.225056415159502
$ws.Range("E5").Value = 9.692575405310716
$ws.Range("F5").Value = 41.8276175720948
$ws.Range("H5").Value = 7.344005520526261
$ws.Range("L5").Value = 10.05446503871265
$ws.Range("M5").Value = 18.32233977751599
$ws.Range("N5").Value = 21.41832211871932
$ws.Range("B6").Value = 22.43760839615973
$ws.Range("C6").Value = 7.936606157333362
$ws.Range("D6").Value = 8.226145633695998
$ws.Range("E6").Value = 9.691801075611279
$ws.Range("F6").Value = 41.8162775593719
$ws.Range("H6").Value = 7.344005520526261
$ws.Range("L6").Value = 10.05520108717855
$ws.Range("M6").Value = 18.32000781733085
$ws.Range("N6").Value = 21.42072234125568
$ws.Range("B7").Value = 22.55092262584592
$ws.Range("C7").Value = 8.127339245173227
$ws.Range("D7").Value = 8.218676965071195
$ws.Range("E7").Value = 9.697152835420395
$ws.Range("F7").Value = 41.8959961258626
$ws.Range("H7").Value = 7.344005520526261
$ws.Range("L7").Value = 10.05022544400875
$ws.Range("M7").Value = 18.3365865196763
$ws.Range("N7").Value = 21.40423632434928
$ws.Range("B8").Value = 23.06032937887678
$ws.Range("C8").Value = 8.918847913961097
$ws.Range("D8").Value = 8.187991328453943
$ws.Range("E8").Value = 9.720203551008076
$ws.Range("F8").Value = 42.27232972341312
$ws.Range("H8").Value = 7.344005520526261
$ws.Range("L8").Value = 10.03156949000335
$ws.Range("M8").Value = 18.41939007785528
$ws.Range("N8").Value = 21.33579884677151
$ws.Range("B9").Value = 24.08372626987843
$ws.Range("C9").Value = 10.29183033728992
$ws.Range("D9").Value = 8.136011402257251
$ws.Range("E9").Value = 9.76339481053034
$ws.Range("F9").Value = 43.09479510621736
$ws.Range("H9").Value = 7.344005520526261
$ws.Range("L9").Value = 10.00685491608031
$ws.Range("M9").Value = 18.6157304367277
$ws.Range("N9").Value = 21.21711358252189
$ws.Range("B10").Value = 24.84233232133093
$ws.Range("C10").Value = 11.19242253867685
$ws.Range("D10").Value = 8.102844657830788
$ws.Range("E10").Value = 9.79392349220857
$ws.Range("F10").Value = 43.74674981790645
$ws.Range("H10").Value = 7.344005520526261
$ws.Range("L10").Value = 9.995938068887703
$ws.Range("M10").Value = 18.77994718712089
$ws.Range("N10").Value = 21.13939850484956
$ws.Range("B11").Value = 25.18721969605732
$ws.Range("C11").Value = 11.57888310752609
$ws.Range("D11").Value = 8.088853209887965
$ws.Range("E11").Value = 9.807555181247986
$ws.Range("F11").Value = 44.05279782175253
$ws.Range("H11").Value = 7.344005520526261
$ws.Range("L11").Value = 9.992544050040459
$ws.Range("M11").Value = 18.85880601053776
$ws.Range("N11").Value = 21.10610828936532
$ws.Range("B12").Value = 25.31765570172889
$ws.Range("C12").Value = 11.72190508491963
$ws.Range("D12").Value = 8.083713233337054
$ws.Range("E12").Value = 9.812680677661607
$ws.Range("F12").Value = 44.16997138394821
$ws.Range("H12").Value = 7.344005520526261
$ws.Range("L12").Value = 9.991484767435765
$ws.Range("M12").Value = 18.88924816669322
$ws.Range("N12").Value = 21.09379934994232
$ws.Range("B13").Value = 25.28957342888545
$ws.Range("C13").Value = 11.69125022546731
$ws.Range("D13").Value = 8.084813169238242
$ws.Range("E13").Value = 9.811578430485779
$ws.Range("F13").Value = 44.14468037379734
$ws.Range("H13").Value = 7.344005520526261
$ws.Range("L13").Value = 9.991702855708457
$ws.Range("M13").Value = 18.88266639715575
$ws.Range("N13").Value = 21.09643706397651
$ws.Range("B14").Value = 25.19795481599297
$ws.Range("C14").Value = 11.59071609541974
$ws.Range("D14").Value = 8.088427164997846
$ws.Range("E14").Value = 9.807977589408051
$ws.Range("F14").Value = 44.06241259076016
$ws.Range("H14").Value = 7.344005520526261
$ws.Range("L14").Value = 9.992452374515864
$ws.Range("M14").Value = 18.86129898460884
$ws.Range("N14").Value = 21.10508966206527
$ws.Range("B15").Value = 25.14181030406568
$ws.Range("C15").Value = 11.52870389615456
$ws.Range("D15").Value = 8.090661473962104
$ws.Range("E15").Value = 9.80576722147471
$ws.Range("F15").Value = 44.0121853897727
$ws.Range("H15").Value = 7.344005520526261
$ws.Range("L15").Value = 9.992940898642356
$ws.Range("M15").Value = 18.84828582577384
$ws.Range("N15").Value = 21.11042837304133
$ws.Range("B16").Value = 24.81977724399716
$ws.Range("C16").Value = 11.16670058903792
$ws.Range("D16").Value = 8.103781145716448
$ws.Range("E16").Value = 9.793027502490837
$ws.Range("F16").Value = 43.72693248279056
$ws.Range("H16").Value = 7.344005520526261
$ws.Range("L16").Value = 9.996191502963603
$ws.Range("M16").Value = 18.77487572480193
$ws.Range("N16").Value = 21.14161566743206
$ws.Range("B17").Value = 24.62206928776162
$ws.Range("C17").Value = 10.93868745188481
$ws.Range("D17").Value = 8.112110867959634
$ws.Range("E17").Value = 9.785146737974443
$ws.Range("F17").Value = 43.55430721945754
$ws.Range("H17").Value = 7.344005520526261
$ws.Range("L17").Value = 9.998588240377169
$ws.Range("M17").Value = 18.73089323878386
$ws.Range("N17").Value = 21.16127688529575
$ws.Range("B18").Value = 24.50834127544064
$ws.Range("C18").Value = 10.8053533941606
$ws.Range("D18").Value = 8.117005050502724
$ws.Range("E18").Value = 9.780589821726252
$ws.Range("F18").Value = 43.45591424834608
$ws.Range("H18").Value = 7.344005520526261
$ws.Range("L18").Value = 10.00011476515286
$ws.Range("M18").Value = 18.70598796012189
$ws.Range("N18").Value = 21.17277960944947
$ws.Range("B19").Value = 24.46983715166859
$ws.Range("C19").Value = 10.75983264275759
$ws.Range("D19").Value = 8.118679833027713
$ws.Range("E19").Value = 9.779042778975072
$ws.Range("F19").Value = 43.4227565101237
$ws.Range("H19").Value = 7.344005520526261
$ws.Range("L19").Value = 10.00065703904533
$ws.Range("M19").Value = 18.6976233539963
$ws.Range("N19").Value = 21.17670755170897
$ws.Range("B20").Value = 24.64311784085452
$ws.Range("C20").Value = 10.96318605421997
$ws.Range("D20").Value = 8.111213474971359
$ws.Range("E20").Value = 9.785988151974037
$ws.Range("F20").Value = 43.57259127322114
$ws.Range("H20").Value = 7.344005520526261
$ws.Range("L20").Value = 9.998317788320341
$ws.Range("M20").Value = 18.73553477358592
$ws.Range("N20").Value = 21.15916382018943
$ws.Range("B21").Value = 25.22487096472711
$ws.Range("C21").Value = 11.62033543342609
$ws.Range("D21").Value = 8.08736134593047
$ws.Range("E21").Value = 9.809036232870152
$ws.Range("F21").Value = 44.08654256538455
$ws.Range("H21").Value = 7.344005520526261
$ws.Range("L21").Value = 9.992226091390433
$ws.Range("M21").Value = 18.86755951480227
$ws.Range("N21").Value = 21.10254010853203
$ws.Range("B22").Value = 25.60405759415331
$ws.Range("C22").Value = 12.03046398361793
$ws.Range("D22").Value = 8.072695496652806
$ws.Range("E22").Value = 9.823886700131142
$ws.Range("F22").Value = 44.42985948230282
$ws.Range("H22").Value = 7.344005520526261
$ws.Range("L22").Value = 9.989561781480305
$ws.Range("M22").Value = 18.95721733290635
$ws.Range("N22").Value = 21.0672667422397
$ws.Range("B23").Value = 25.40181547170266
$ws.Range("C23").Value = 11.81333606067724
$ws.Range("D23").Value = 8.080438274486056
$ws.Range("E23").Value = 9.815980104165812
$ws.Range("F23").Value = 44.24597359325298
$ws.Range("H23").Value = 7.344005520526261
$ws.Range("L23").Value = 9.990863323108073
$ws.Range("M23").Value = 18.90906279063191
$ws.Range("N23").Value = 21.0859339379825
$ws.Range("B24").Value = 24.63360197610892
$ws.Range("C24").Value = 10.95211723403113
$ws.Range("D24").Value = 8.111618858359058
$ws.Range("E24").Value = 9.785607830207473
$ws.Range("F24").Value = 43.56432239380567
$ws.Range("H24").Value = 7.344005520526261
$ws.Range("L24").Value = 9.998439596805278
$ws.Range("M24").Value = 18.73343514815542
$ws.Range("N24").Value = 21.16011851642699
$ws.Range("B25").Value = 23.8051071963515
$ws.Range("C25").Value = 9.939472193788134
$ws.Range("D25").Value = 8.149193326718756
$ws.Range("E25").Value = 9.751923445189288
$ws.Range("F25").Value = 42.8636283525552
$ws.Range("H25").Value = 7.344005520526261
$ws.Range("L25").Value = 10.01226940659432
$ws.Range("M25").Value = 18.55905051383602
$ws.Range("N25").Value = 21.24755745484794
